$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: C1 text changes, D1/E1 are new headers copied from C1 style
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("C1").Copy($ws.Range("E1"))
$ws.Range("C1").Value = "Frecuencia del primer armonico"
$ws.Range("D1").Value = "Frecuencia del segundo armonico"
$ws.Range("E1").Value = "Frecuencia tercer armonico"

# Update data rows 2-48 for columns C (revised), D (new), E (new)
$ws.Range("C2").Value = 160.9074677029721
$ws.Range("D2").Value = 323.9155550887517
$ws.Range("E2").Value = 486.9236424745304
$ws.Range("C3").Value = 313.5672986419258
$ws.Range("D3").Value = 156.5147236789026
$ws.Range("E3").Value = 468.4684684684689
$ws.Range("C4").Value = 186.4502736365348
$ws.Range("D4").Value = 371.3908284581994
$ws.Range("E4").Value = 557.8411020947347
$ws.Range("C5").Value = 271.844660194175
$ws.Range("D5").Value = 136.6990291262136
$ws.Range("E5").Value = 406.9902912621355
$ws.Range("C6").Value = 303.0653556969346
$ws.Range("D6").Value = 152.358919276212
$ws.Range("E6").Value = 455.7547715442452
$ws.Range("C7").Value = 148.8622669596652
$ws.Range("D7").Value = 297.4409867441691
$ws.Range("E7").Value = 446.0197065286729
$ws.Range("C8").Value = 209.9517111064451
$ws.Range("D8").Value = 417.384001679613
$ws.Range("E8").Value = 619.7774511862272
$ws.Range("C9").Value = 371.350720124562
$ws.Range("D9").Value = 185.2861035422343
$ws.Range("E9").Value = 2407.162319968859
$ws.Range("C10").Value = 308.6115992970126
$ws.Range("D10").Value = 154.6572934973638
$ws.Range("E10").Value = 461.1599297012308
$ws.Range("C11").Value = 379.6507213363711
$ws.Range("D11").Value = 187.6993166287016
$ws.Range("E11").Value = 2291.268033409263
$ws.Range("C12").Value = 303.5591274397248
$ws.Range("D12").Value = 606.1997703788747
$ws.Range("E12").Value = 250.7462686567164
$ws.Range("C13").Value = 314.8374480567099
$ws.Range("D13").Value = 620.8750916646295
$ws.Range("E13").Value = 3090.686873625031
$ws.Range("C14").Value = 401.6001561127914
$ws.Range("D14").Value = 201.7757830032197
$ws.Range("E14").Value = 602.2050931798226
$ws.Range("C15").Value = 171.1229946524063
$ws.Range("D15").Value = 342.2459893048126
$ws.Range("E15").Value = 513.3689839572194
$ws.Range("C16").Value = 189.0133333333333
$ws.Range("D16").Value = 375.8933333333334
$ws.Range("E16").Value = 2648.320000000001
$ws.Range("C17").Value = 331.593490942585
$ws.Range("D17").Value = 165.7967454712925
$ws.Range("E17").Value = 490.0214921707093
$ws.Range("C18").Value = 391.5139826422374
$ws.Range("D18").Value = 195.4355512696884
$ws.Range("E18").Value = 586.9495339119258
$ws.Range("C19").Value = 360.3771823359166
$ws.Range("D19").Value = 184.8566893847442
$ws.Range("E19").Value = 544.1135281486322
$ws.Range("C20").Value = 323.8721242904094
$ws.Range("D20").Value = 162.931978886565
$ws.Range("E20").Value = 485.609003087342
$ws.Range("C21").Value = 312.4042879019908
$ws.Range("D21").Value = 155.7646029315247
$ws.Range("E21").Value = 467.2938087945749
$ws.Range("C22").Value = 355.0190029063269
$ws.Range("D22").Value = 178.403755868545
$ws.Range("E22").Value = 538.7882852671587
$ws.Range("C23").Value = 312.6536230459146
$ws.Range("D23").Value = 157.703274014355
$ws.Range("E23").Value = 467.6039720774752
$ws.Range("C24").Value = 165.0294695481334
$ws.Range("D24").Value = 331.3686967910935
$ws.Range("E24").Value = 497.7079240340536
$ws.Range("C25").Value = 264.368869740581
$ws.Range("D25").Value = 524.7043725364383
$ws.Range("E25").Value = 2908.057567146393
$ws.Range("C26").Value = 246.3994447336454
$ws.Range("D26").Value = 493.0302504482624
$ws.Range("E26").Value = 3678.408236450923
$ws.Range("C27").Value = 279.5876288659792
$ws.Range("D27").Value = 560
$ws.Range("E27").Value = 2814.845360824742
$ws.Range("C28").Value = 261.8162506638341
$ws.Range("D28").Value = 525.225703664365
$ws.Range("E28").Value = 2877.854487519915
$ws.Range("C29").Value = 491.2612187057157
$ws.Range("D29").Value = 244.685876239962
$ws.Range("E29").Value = 735.9470949456781
$ws.Range("C30").Value = 239.8433675966717
$ws.Range("D30").Value = 479.6867351933433
$ws.Range("E30").Value = 719.530102790015
$ws.Range("C31").Value = 258.4615384615381
$ws.Range("D31").Value = 517.3474801061006
$ws.Range("E31").Value = 2900.79575596817
$ws.Range("C32").Value = 203.9101041476333
$ws.Range("D32").Value = 399.0498812351543
$ws.Range("E32").Value = 2188.196601498264
$ws.Range("C33").Value = 220.9462191670036
$ws.Range("D33").Value = 440.2749696724622
$ws.Range("E33").Value = 2871.65386170643
$ws.Range("C34").Value = 235.7335546793938
$ws.Range("D34").Value = 472.2971570865325
$ws.Range("E34").Value = 2851.213944801826
$ws.Range("C35").Value = 292.1088578333511
$ws.Range("D35").Value = 586.3192182410421
$ws.Range("E35").Value = 183.6713249973732
$ws.Range("C36").Value = 295.2351292432268
$ws.Range("D36").Value = 591.7159763313612
$ws.Range("E36").Value = 2964.39323159971
$ws.Range("C37").Value = 247.3241192975165
$ws.Range("D37").Value = 499.2206172711212
$ws.Range("E37").Value = 2747.168242751741
$ws.Range("C38").Value = 200.472316977171
$ws.Range("D38").Value = 428.2340593020208
$ws.Range("E38").Value = 595.1193912358958
$ws.Range("C39").Value = 278.6592413575709
$ws.Range("D39").Value = 555.2169801408008
$ws.Range("E39").Value = 444.2576442156142
$ws.Range("C40").Value = 290.8685140412872
$ws.Range("D40").Value = 578.7613911102844
$ws.Range("E40").Value = 427.0039055235266
$ws.Range("C41").Value = 282.0211515863693
$ws.Range("D41").Value = 559.3419506462988
$ws.Range("E41").Value = 2820.99490795143
$ws.Range("C42").Value = 323.0109575988572
$ws.Range("D42").Value = 161.0290614578371
$ws.Range("E42").Value = 484.992853739876
$ws.Range("C43").Value = 289.0303257207033
$ws.Range("D43").Value = 433.5454885810559
$ws.Range("E43").Value = 144.5151628603517
$ws.Range("C44").Value = 201.6677357280314
$ws.Range("D44").Value = 401.7960230917261
$ws.Range("E44").Value = 602.9506093649775
$ws.Range("C45").Value = 194.4191438271091
$ws.Range("D45").Value = 388.5070795727415
$ws.Range("E45").Value = 581.9325991554192
$ws.Range("C46").Value = 208.4765177548679
$ws.Range("D46").Value = 416.9530355097359
$ws.Range("E46").Value = 626.5750286368839
$ws.Range("C47").Value = 402.5996838222381
$ws.Range("D47").Value = 200.9485332864915
$ws.Range("E47").Value = 604.250834357983
$ws.Range("C48").Value = 223.8885913229778
$ws.Range("D48").Value = 448.8484199250133
$ws.Range("E48").Value = 2687.198714515265
